# Auto-generated edit script applying numeric cell updates
# as described by the source diff across all 8 sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 873.2308
$ws.Range("I41").Value = 803.6667
$ws.Range("K41").Value = 803.6667
$ws.Range("M41").Value = -363.6667
$ws.Range("H64").Value = 9636.362999999999
$ws.Range("J64").Value = 9636.362999999999
$ws.Range("L64").Value = 9636.362999999999
$ws.Range("N64").Value = -10132.363
$ws.Range("H67").Value = 9636.362999999999
$ws.Range("J67").Value = 9636.362999999999
$ws.Range("L67").Value = 9636.362999999999
$ws.Range("N67").Value = -11352.363
$ws.Range("H70").Value = 92314850
$ws.Range("I70").Value = 50004500
$ws.Range("K70").Value = 150013500
$ws.Range("M70").Value = -150013230
$ws.Range("H73").Value = 92314850
$ws.Range("I73").Value = 50004500
$ws.Range("K73").Value = 150013500
$ws.Range("M73").Value = -150012564
$ws.Range("H100").Value = 5405.0527
$ws.Range("I100").Value = 4582.8335
$ws.Range("K100").Value = 4582.8335
$ws.Range("M100").Value = -4041.8335
$ws.Range("H111").Value = 720
$ws.Range("I111").Value = 700
$ws.Range("K111").Value = 2100
$ws.Range("M111").Value = 967

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1982
$ws.Range("J2").Value = 3313.889
$ws.Range("L2").Value = 3313.889
$ws.Range("N2").Value = -3539.889
$ws.Range("H5").Value = 67.666664
$ws.Range("J5").Value = 70.666664
$ws.Range("L5").Value = 70.666664
$ws.Range("N5").Value = -294.666664
$ws.Range("H11").Value = 5251251
$ws.Range("I11").Value = 7000001
$ws.Range("K11").Value = 7000001
$ws.Range("M11").Value = -6999857
$ws.Range("H32").Value = 3988
$ws.Range("I32").Value = 3988
$ws.Range("K32").Value = 3988
$ws.Range("M32").Value = -3701
$ws.Range("H45").Value = 983
$ws.Range("I45").Value = 983
$ws.Range("K45").Value = 983
$ws.Range("M45").Value = -606
$ws.Range("H61").Value = 3750.9
$ws.Range("I61").Value = 3776.7932
$ws.Range("K61").Value = 3776.7932
$ws.Range("M61").Value = -3564.7932
$ws.Range("H63").Value = 6284.4287
$ws.Range("I63").Value = 499.5
$ws.Range("K63").Value = 499.5
$ws.Range("M63").Value = 186.5
$ws.Range("H66").Value = 6284.4287
$ws.Range("I66").Value = 499.5
$ws.Range("K66").Value = 2497.5
$ws.Range("M66").Value = 934.5
$ws.Range("H102").Value = 3622.2593
$ws.Range("I102").Value = 2172.7727
$ws.Range("K102").Value = 2172.7727
$ws.Range("M102").Value = -550.7727
$ws.Range("H110").Value = 662.5
$ws.Range("I110").Value = 662.5
$ws.Range("K110").Value = 662.5
$ws.Range("M110").Value = 1382.5
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents() | Out-Null
$ws.Range("H116").Value = 1982
$ws.Range("J116").Value = 3313.889
$ws.Range("L116").Value = 3313.889
$ws.Range("N116").Value = -7901.889
$ws.Range("H136").Value = 3750.9
$ws.Range("I136").Value = 3776.7932
$ws.Range("K136").Value = 11330.3796
$ws.Range("M136").Value = -8780.3796

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1982
$ws.Range("J3").Value = 3313.889
$ws.Range("L3").Value = 3313.889
$ws.Range("N3").Value = -3541.889
$ws.Range("H4").Value = 67.666664
$ws.Range("J4").Value = 70.666664
$ws.Range("L4").Value = 70.666664
$ws.Range("N4").Value = -300.666664
$ws.Range("H99").Value = 2971.625
$ws.Range("I99").Value = 1687.3334
$ws.Range("J99").Value = 4622.857
$ws.Range("K99").Value = 1687.3334
$ws.Range("L99").Value = 4622.857
$ws.Range("M99").Value = -189.3334
$ws.Range("N99").Value = -7618.857
$ws.Range("H134").Value = 3924.9443
$ws.Range("J134").Value = 3199.8333
$ws.Range("L134").Value = 9599.499899999999
$ws.Range("N134").Value = -14669.4999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1268.1111
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents() | Out-Null
$ws.Range("H99").Value = 3016.6667
$ws.Range("I99").Value = 3016.6667
$ws.Range("K99").Value = 3016.6667
$ws.Range("M99").Value = -1518.6667
$ws.Range("H107").Value = 1289.2307
$ws.Range("I107").Value = 1289.2307
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1289.2307
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 630.7692999999999
$ws.Range("N107").ClearContents() | Out-Null
$ws.Range("H122").Value = 4269.846
$ws.Range("I122").Value = 4367.778
$ws.Range("K122").Value = 13103.334
$ws.Range("M122").Value = -10653.334
$ws.Range("H126").Value = 3016.6667
$ws.Range("I126").Value = 3016.6667
$ws.Range("K126").Value = 9050.000100000001
$ws.Range("M126").Value = -6580.000100000001
$ws.Range("H132").Value = 1880.2963
$ws.Range("I132").Value = 1832
$ws.Range("J132").Value = 2484
$ws.Range("K132").Value = 5496
$ws.Range("L132").Value = 7452
$ws.Range("M132").Value = -2966
$ws.Range("N132").Value = -12512
$ws.Range("H136").Value = 1268.1111
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents() | Out-Null

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 4500
$ws.Range("I57").Value = 2000
$ws.Range("J57").Value = 5125
$ws.Range("K57").Value = 6000
$ws.Range("L57").Value = 15375
$ws.Range("M57").Value = -5441
$ws.Range("N57").Value = -16493
$ws.Range("H64").Value = 999.6
$ws.Range("J64").Value = 1004.6667
$ws.Range("L64").Value = 3014.0001
$ws.Range("N64").Value = -3554.0001
$ws.Range("H67").Value = 999.6
$ws.Range("J67").Value = 1004.6667
$ws.Range("L67").Value = 3014.0001
$ws.Range("N67").Value = -4886.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 13897154
$ws.Range("I70").Value = 83336776
$ws.Range("K70").Value = 83336776
$ws.Range("M70").Value = -83336506
$ws.Range("H73").Value = 13897154
$ws.Range("I73").Value = 83336776
$ws.Range("K73").Value = 83336776
$ws.Range("M73").Value = -83335840
$ws.Range("H107").Value = 378.42856
$ws.Range("J107").Value = 601
$ws.Range("L107").Value = 601
$ws.Range("N107").Value = -4441

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 7715
$ws.Range("I18").Value = 7715
$ws.Range("K18").Value = 7715
$ws.Range("M18").Value = -7543
$ws.Range("H46").Value = 3413.4736
$ws.Range("I46").Value = 791.6667
$ws.Range("K46").Value = 791.6667
$ws.Range("M46").Value = -603.6667
$ws.Range("H136").Value = 2634.0625
$ws.Range("I136").Value = 2206.7144
$ws.Range("K136").Value = 6620.1432
$ws.Range("M136").Value = -4070.1432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H37").Value = 25000
$ws.Range("J37").Value = 25000
$ws.Range("L37").Value = 25000
$ws.Range("N37").Value = -25406
$ws.Range("H62").Value = 10875
$ws.Range("I62").Value = 7500
$ws.Range("J62").Value = 12000
$ws.Range("K62").Value = 7500
$ws.Range("L62").Value = 12000
$ws.Range("M62").Value = -6876
$ws.Range("N62").Value = -13248
$ws.Range("H65").Value = 10875
$ws.Range("I65").Value = 7500
$ws.Range("J65").Value = 12000
$ws.Range("K65").Value = 37500
$ws.Range("L65").Value = 60000
$ws.Range("M65").Value = -34380
$ws.Range("N65").Value = -66240
